$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "143×9=1287" "900×3=2700"
Replace-Text "399×7=2793" "163×7=1141"
Replace-Text "142×5=710" "199×4=796"
Replace-Text "805×8=6440" "433×3=1299"
Replace-Text "885×2=1770" "325×5=1625"
Replace-Text "588×6=3528" "843×9=7587"
Replace-Text "183×3=549" "195×9=1755"
Replace-Text "414×9=3726" "309×2=618"
Replace-Text "922×4=3688" "994×3=2982"
Replace-Text "227×5=1135" "989×4=3956"
Replace-Text "118×5=590" "561×4=2244"
Replace-Text "847×9=7623" "978×5=4890"
Replace-Text "953×7=6671" "912×7=6384"
Replace-Text "218×5=1090" "264×7=1848"
Replace-Text "531×2=1062" "135×8=1080"
Replace-Text "304×5=1520" "385×7=2695"
Replace-Text "866×8=6928" "916×5=4580"
Replace-Text "518×9=4662" "934×6=5604"
Replace-Text "976×2=1952" "510×7=3570"
Replace-Text "638×3=1914" "655×8=5240"
Replace-Text "128×8=1024" "132×5=660"
Replace-Text "929×6=5574" "243×2=486"
Replace-Text "101×9=909" "291×6=1746"
Replace-Text "807×2=1614" "723×8=5784"
Replace-Text "325×9=2925" "267×4=1068"
